$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$ws.Name = "L6"
$ws.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws = $wb.Worksheets.Item("L6")

$ws.Range("B1").Value = 'Form'
$ws.Range("C1").Value = 'Goals scored'
$ws.Range("D1").Value = 'Goals conceded'
$ws.Range("E1").Value = 'Total Goals'

$bData = @(
    'Barrow,W L L L W L',
    'Bolton,W W L W W L',
    'Bradford,L L L L L D',
    'Cambridge,D L W W L L',
    'Carlisle,D D L D D W',
    'Cheltenham,W L W W D L',
    'Colchester,D L W W L W',
    'Crawley Town,D W W L D D',
    'Exeter,L D D D W W',
    'Forest Green,L W D D L W',
    'Grimsby,L W L W L W',
    'Harrogate,L D W L D W',
    'Leyton Orient,L D W L L L',
    'Mansfield,D D W W L W',
    'Morecambe,L W W W L W',
    'Newport County,D L W D W W',
    'Oldham,W W L W L L',
    'Port Vale,W W D W W L',
    'Salford,L D W W W L',
    'Scunthorpe,L L L L L D',
    'Southend,D D D L W W',
    'Stevenage,D L L L W D',
    'Tranmere,D L D W W L',
    'Walsall,W W L L W L'
)
$cData = @(
    'Barrow,2 0 0 0 2 1',
    'Bolton,2 1 1 1 1 1',
    'Bradford,0 1 0 1 0 0',
    'Cambridge,1 1 1 4 0 4',
    'Carlisle,0 0 0 1 1 3',
    'Cheltenham,1 0 1 1 1 0',
    'Colchester,1 2 2 2 0 1',
    'Crawley Town,0 2 1 0 0 3',
    'Exeter,1 0 0 0 3 2',
    'Forest Green,1 3 0 0 0 2',
    'Grimsby,0 2 0 2 2 1',
    'Harrogate,1 2 2 0 1 5',
    'Leyton Orient,0 2 2 2 1 2',
    'Mansfield,0 1 1 3 0 4',
    'Morecambe,0 4 4 3 0 2',
    'Newport County,0 0 2 0 4 1',
    'Oldham,4 5 3 3 1 1',
    'Port Vale,2 1 0 2 2 0',
    'Salford,0 0 2 2 1 0',
    'Scunthorpe,1 2 0 0 0 0',
    'Southend,1 0 0 0 2 2',
    'Stevenage,1 1 0 0 1 3',
    'Tranmere,0 0 0 1 1 1',
    'Walsall,2 1 1 0 2 0'
)
$dData = @(
    'Barrow,1 2 2 1 0 2',
    'Bolton,1 0 2 0 0 2',
    'Bradford,2 2 1 2 1 0',
    'Cambridge,1 4 0 2 1 5',
    'Carlisle,0 0 1 1 1 2',
    'Cheltenham,0 1 0 0 1 1',
    'Colchester,1 5 1 0 1 0',
    'Crawley Town,0 0 0 2 0 3',
    'Exeter,2 0 0 0 2 1',
    'Forest Green,2 2 0 0 2 1',
    'Grimsby,1 1 3 1 3 0',
    'Harrogate,2 2 1 3 1 4',
    'Leyton Orient,1 2 0 4 2 3',
    'Mansfield,0 1 0 0 2 1',
    'Morecambe,1 1 3 0 1 0',
    'Newport County,0 1 0 0 0 0',
    'Oldham,1 2 4 0 2 4',
    'Port Vale,0 0 0 0 1 1',
    'Salford,1 0 0 0 0 1',
    'Scunthorpe,4 3 3 2 4 0',
    'Southend,1 0 0 2 1 1',
    'Stevenage,1 2 1 1 0 3',
    'Tranmere,0 1 0 0 0 2',
    'Walsall,1 0 2 2 0 2'
)
$eData = @(
    'Barrow,3 2 2 1 2 3',
    'Bolton,3 1 3 1 1 3',
    'Bradford,2 3 1 3 1 0',
    'Cambridge,2 5 1 6 1 9',
    'Carlisle,0 0 1 2 2 5',
    'Cheltenham,1 1 1 1 2 1',
    'Colchester,2 7 3 2 1 1',
    'Crawley Town,0 2 1 2 0 6',
    'Exeter,3 0 0 0 5 3',
    'Forest Green,3 5 0 0 2 3',
    'Grimsby,1 3 3 3 5 1',
    'Harrogate,3 4 3 3 2 9',
    'Leyton Orient,1 4 2 6 3 5',
    'Mansfield,0 2 1 3 2 5',
    'Morecambe,1 5 7 3 1 2',
    'Newport County,0 1 2 0 4 1',
    'Oldham,5 7 7 3 3 5',
    'Port Vale,2 1 0 2 3 1',
    'Salford,1 0 2 2 1 1',
    'Scunthorpe,5 5 3 2 4 0',
    'Southend,2 0 0 2 3 3',
    'Stevenage,2 3 1 1 1 6',
    'Tranmere,0 1 0 1 1 3',
    'Walsall,3 1 3 2 2 2'
)

for ($i = 0; $i -lt 24; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = [string]($i + 1)
    $ws.Cells.Item($row, 2).Value = $bData[$i]
    $ws.Cells.Item($row, 3).Value = $cData[$i]
    $ws.Cells.Item($row, 4).Value = $dData[$i]
    $ws.Cells.Item($row, 5).Value = $eData[$i]
}
